# Added new Matrix problem ("Number of Islands") as a new row at the
# bottom of the Notes table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row's contents. Column D (the link) is entered before
# column C (the algorithm write-up) to reproduce the shared-string
# allocation order of the authored workbook.
$ws.Cells.Item(21, 1).Value = "Number of Islands"
$ws.Cells.Item(21, 2).Value = "Return num of islands made from adjacent '1' "
$ws.Cells.Item(21, 4).Value = "https://leetcode.com/problems/number-of-islands/"
$ws.Cells.Item(21, 3).Value = "Iterate over grid using 2 for loops. If current is '1' increase numOfIslands, call recursive DFS function to mark all as visited(change it to '0')"

# Hyperlink the new URL cell to its own text.
$ws.Hyperlinks.Add($ws.Cells.Item(21, 4), "https://leetcode.com/problems/number-of-islands/") | Out-Null

# Copy the formatting of the last existing row (20) down into the new
# row (21) (after adding the hyperlink) so every cell - including the
# link cell - ends up on the same styles already used by the rest of
# the table (name/desc/algo/link columns).
$ws.Range("A20:D20").Copy() | Out-Null
$ws.Range("A21:D21").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

# Move the active selection to C22, matching where the cursor ended up
# after the author finished typing the new row.
$ws.Range("C22").Select() | Out-Null
